$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments -------------------------------------------------
# NOTE: the engine stores OOXML <col width> as ColumnWidth + 0.83 (standard
# Excel character-width padding for the default font), so back the target
# stored width off by 0.83 before assigning ColumnWidth.
$ws.Columns.Item(2).ColumnWidth = 38 - 0.83
$ws.Columns.Item(4).ColumnWidth = 28 - 0.83
$ws.Columns.Item(8).ColumnWidth = 22 - 0.83

# --- Refresh the scraped listing data (rows 2-5) ------------------------------
$ws.Range("A2").Value = "2026-01-31 06:35:46"
$ws.Range("B2").Value = "【AI活用】業務改善DXツール開発エンジニア募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5482904"
$ws.Range("G2").Value = 423
$ws.Range("H2").Value = "🔥AI,Ai ◆ツール,開発 ◇業務改善"

$ws.Range("A3").Value = "2026-01-31 06:35:46"
$ws.Range("B3").Value = "【急募】NTT線路情報開示システムの自動化依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5482939"
$ws.Range("G3").Value = 98
$ws.Range("H3").Value = "◆自動化"

$ws.Range("A4").Value = "2026-01-31 06:35:46"
$ws.Range("B4").Value = "【Excel管理表】計算式保持の修正・最適化依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5482932"
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = "◇管理"

$ws.Range("A5").Value = "2026-01-31 06:35:46"
$ws.Range("B5").Value = "Power Automate(またはGAS)での予約サイト連携フロー構築"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5482835"
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = "◇サイト"

# --- Drop the old rows that fell out of the refreshed top-4 listing ----------
$ws.Range("A6:H13").Delete()

# --- Rebuild the hyperlinks so only F2:F5 keep live links ---------------------
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5482904") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5482939") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5482932") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5482835") | Out-Null
